# Add a new code-snippets row (row 6) for the "Suppress Plots and Results" Rmd snippet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# title, src, url, description, tags, All, R
$ws.Range("A6").Value = "Hide certain plots and results in rendered Rmd"
$ws.Range("B6").Value = "images/arseny-togulev-upnf6XRkWho-unsplash.jpg"
$ws.Range("C6").Value = "https://sciencificity.github.io/rmd-hide-info/"
$ws.Range("D6").Value = "Use code chunk options to hide results and plots"
$ws.Range("E6").Value = "R; Suppress Plots and Results in report"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1

# Turn the url cell into a hyperlink, same as the other rows.
$ws.Hyperlinks.Add($ws.Range("C6"), "https://sciencificity.github.io/rmd-hide-info/")

# Re-apply the "Hyperlink" cell style (matching C2/C4/C5) since Hyperlinks.Add
# creates its own style entry; copying the format from an existing linked cell
# keeps the same style index used elsewhere in the sheet.
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Matches the saved selection recorded in the workbook.
$ws.Range("C18").Select()

$wb.Save()
